# MLB Stats 2019.xlsx - "Add files via upload" edit
#
# Summary of content changes (per the OOXML diff):
#  1. Sheet "2019 League Hitting" (index 2): the batting-average column
#     header in R1 is renamed from "BA" to "AVG" (the old "BA" shared
#     string is dropped since nothing else referenced it).
#  2. Sheet "2019 League Hitting": a new "wRC+" column is appended in
#     column AE (right after WAR in AD), with one integer value per
#     player row (2-31) and an AVERAGE formula in row 32, matching the
#     existing WAR column's look (same number formats / styles).
#  3. The workbook-level revisionPtr / window-size bits and the sheet's
#     view scroll position are incidental Excel-generated noise from the
#     save and are not reproducible content edits, so they're left alone
#     except for the selected cell, which we do set to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2019 League Hitting")
$ws.Activate()

# --- 1. Rename the "BA" header to "AVG" ---
$ws.Range("R1").Value = "AVG"

# --- 2. Add the new "wRC+" column in AE ---
$ws.Range("AE1").Value = "wRC+"

$wrcPlus = @(93,101,88,106,100,92,85,97,87,77,124,84,96,110,79,97,116,105,117,107,90,92,89,99,82,95,104,88,92,104)
for ($i = 0; $i -lt $wrcPlus.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 31).Value = $wrcPlus[$i]
}

$ws.Range("AE32").Formula = "=AVERAGE(AE2:AE31)"

# Match the WAR column's formatting (header style, data style, average-row
# style) for the new wRC+ column.
$ws.Range("AD1:AD32").Copy()
$ws.Range("AE1:AE32").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Restore the selection the author left the sheet on ---
$ws.Range("AC35").Select()
